# Auto-generated edit script: updates market-price columns (H-N) across all 8 Sheets
# Source: scheduled runner price refresh (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 102.90909
$ws.Range("I11").Value = 102.90909
$ws.Range("K11").Value = 102.90909
$ws.Range("M11").Value = 37.09090999999999
# Row 12
$ws.Range("H12").Value = 390.1
$ws.Range("J12").Value = 1000
$ws.Range("L12").Value = 1000
$ws.Range("N12").Value = -1340
# Row 80
$ws.Range("H80").Value = 2679.6904
$ws.Range("I80").Value = 1297.7273
$ws.Range("J80").Value = 3170.0645
$ws.Range("K80").Value = 3893.1819
$ws.Range("L80").Value = 9510.193499999999
$ws.Range("M80").Value = -2895.1819
$ws.Range("N80").Value = -11506.1935
# Row 83
$ws.Range("H83").Value = 2679.6904
$ws.Range("I83").Value = 1297.7273
$ws.Range("J83").Value = 3170.0645
$ws.Range("K83").Value = 11679.5457
$ws.Range("L83").Value = 28530.5805
$ws.Range("M83").Value = -6687.545700000001
$ws.Range("N83").Value = -38514.5805
# Row 112
$ws.Range("H112").Value = 2025
$ws.Range("I112").Value = 2100
$ws.Range("J112").Value = 2000
$ws.Range("K112").Value = 6300
$ws.Range("L112").Value = 6000
$ws.Range("M112").Value = -5192
$ws.Range("N112").Value = -8216
# Row 116
$ws.Range("H116").Value = 5762.75
$ws.Range("I116").Value = 4920.4
$ws.Range("J116").Value = 7166.6665
$ws.Range("K116").Value = 4920.4
$ws.Range("L116").Value = 7166.6665
$ws.Range("M116").Value = -1478.4
$ws.Range("N116").Value = -14050.6665
# Row 132
$ws.Range("H132").Value = 597.08826
$ws.Range("I132").Value = 595.55225
$ws.Range("J132").Value = 700
$ws.Range("K132").Value = 1786.65675
$ws.Range("L132").Value = 2100
$ws.Range("M132").Value = 743.3432500000001
$ws.Range("N132").Value = -7160
# Row 138
$ws.Range("H138").Value = 3775.8276
$ws.Range("J138").Value = 4685.7
$ws.Range("L138").Value = 14057.1
$ws.Range("N138").Value = -24337.1

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12689.862
$ws.Range("I32").Value = 9145.629000000001
$ws.Range("K32").Value = 9145.629000000001
$ws.Range("M32").Value = -8858.629000000001
# Row 61
$ws.Range("I61").Value = 3512.3333
$ws.Range("J61").Value = 6272.1816
$ws.Range("K61").Value = 3512.3333
$ws.Range("L61").Value = 6272.1816
$ws.Range("M61").Value = -3300.3333
$ws.Range("N61").Value = -6696.1816
# Row 74
$ws.Range("H74").Value = 3149.5
$ws.Range("I74").Value = 3000
$ws.Range("J74").Value = 3299
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 3299
$ws.Range("M74").Value = -2126
$ws.Range("N74").Value = -5047
# Row 77
$ws.Range("H77").Value = 3149.5
$ws.Range("I77").Value = 3000
$ws.Range("J77").Value = 3299
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 16495
$ws.Range("M77").Value = -10632
$ws.Range("N77").Value = -25231
# Row 132
$ws.Range("H132").Value = 5141.278
$ws.Range("I132").Value = 5256.4
$ws.Range("J132").Value = 4879.636
$ws.Range("K132").Value = 15769.2
$ws.Range("L132").Value = 14638.908
$ws.Range("M132").Value = -13239.2
$ws.Range("N132").Value = -19698.908
# Row 136
$ws.Range("I136").Value = 3512.3333
$ws.Range("J136").Value = 6272.1816
$ws.Range("K136").Value = 10536.9999
$ws.Range("L136").Value = 18816.5448
$ws.Range("M136").Value = -7986.999899999999
$ws.Range("N136").Value = -23916.5448

$ws = $wb.Worksheets.Item("BSM")
# Row 60
$ws.Range("H60").Value = 84250
$ws.Range("J60").Value = 84250
$ws.Range("L60").Value = 84250
$ws.Range("N60").Value = -85448
# Row 80
$ws.Range("H80").Value = 180.36842
$ws.Range("I80").Value = 127.5
$ws.Range("J80").Value = 186.58824
$ws.Range("K80").Value = 127.5
$ws.Range("L80").Value = 186.58824
$ws.Range("M80").Value = 870.5
$ws.Range("N80").Value = -2182.58824
# Row 83
$ws.Range("H83").Value = 180.36842
$ws.Range("I83").Value = 127.5
$ws.Range("J83").Value = 186.58824
$ws.Range("K83").Value = 637.5
$ws.Range("L83").Value = 932.9412000000001
$ws.Range("M83").Value = 4354.5
$ws.Range("N83").Value = -10916.9412
# Row 102
$ws.Range("H102").Value = 25876.285
$ws.Range("I102").Value = 24789.846
$ws.Range("K102").Value = 24789.846
$ws.Range("M102").Value = -21544.846
# Row 107
$ws.Range("H107").Value = 1197.5
$ws.Range("I107").Value = 1197.5
$ws.Range("K107").Value = 1197.5
$ws.Range("M107").Value = 722.5
# Row 131
$ws.Range("H131").Value = 32270.818
$ws.Range("J131").Value = 32270.818
$ws.Range("L131").Value = 32270.818
$ws.Range("N131").Value = -42350.818
# Row 134
$ws.Range("H134").Value = 4297.4546
$ws.Range("I134").Value = 2696.0715
$ws.Range("K134").Value = 8088.2145
$ws.Range("M134").Value = -5553.2145
# Row 137
$ws.Range("H137").Value = 69988.7
$ws.Range("J137").Value = 69988.7
$ws.Range("L137").Value = 69988.7
$ws.Range("N137").Value = -80188.7
# Row 141
$ws.Range("H141").Value = 59982
$ws.Range("J141").Value = 59964
$ws.Range("L141").Value = 59964
$ws.Range("N141").Value = -70324

$ws = $wb.Worksheets.Item("CRP")
# Row 48
$ws.Range("H48").Value = 45231.668
$ws.Range("J48").Value = 45231.668
$ws.Range("L48").Value = 45231.668
$ws.Range("N48").Value = -46183.668
# Row 58
$ws.Range("H58").Value = 1594.7273
$ws.Range("I58").Value = 942.75
$ws.Range("J58").Value = 3333.3333
$ws.Range("K58").Value = 942.75
$ws.Range("L58").Value = 3333.3333
$ws.Range("M58").Value = -739.75
$ws.Range("N58").Value = -3739.3333
# Row 135
$ws.Range("H135").Value = 69999.35000000001
$ws.Range("J135").Value = 69999.35000000001
$ws.Range("L135").Value = 69999.35000000001
$ws.Range("N135").Value = -80139.35000000001
# Row 136
$ws.Range("H136").Value = 1594.7273
$ws.Range("I136").Value = 942.75
$ws.Range("J136").Value = 3333.3333
$ws.Range("K136").Value = 2828.25
$ws.Range("L136").Value = 9999.999899999999
$ws.Range("M136").Value = -278.25
$ws.Range("N136").Value = -15099.9999

$ws = $wb.Worksheets.Item("CUL")
# Row 140
$ws.Range("H140").Value = 2639.2856
$ws.Range("I140").Value = 1280.5714
$ws.Range("K140").Value = 3841.7142
$ws.Range("M140").Value = 1338.2858

$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 6000
$ws.Range("I9").Value = 6000
$ws.Range("K9").Value = 6000
$ws.Range("M9").Value = -5830
# Row 48
$ws.Range("H48").Value = 30000
$ws.Range("J48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30970
# Row 122
$ws.Range("H122").Value = 3869.4062
$ws.Range("I122").Value = 1530.55
$ws.Range("K122").Value = 4591.65
$ws.Range("M122").Value = -2141.65
# Row 126
$ws.Range("H126").Value = 6653.222
$ws.Range("I126").Value = 1950
$ws.Range("K126").Value = 5850
$ws.Range("M126").Value = -3380
# Row 132
$ws.Range("H132").Value = 4787.5317
$ws.Range("I132").Value = 3903.92
$ws.Range("K132").Value = 11711.76
$ws.Range("M132").Value = -9181.76
# Row 139
$ws.Range("H139").Value = 129666.664
$ws.Range("J139").Value = 129666.664
$ws.Range("L139").Value = 129666.664
$ws.Range("N139").Value = -139946.664

$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 3377
$ws.Range("I9").Value = 5505
$ws.Range("J9").Value = 1249
$ws.Range("K9").Value = 5505
$ws.Range("L9").Value = 1249
$ws.Range("M9").Value = -5281
$ws.Range("N9").Value = -1697
# Row 30
$ws.Range("H30").Value = 5957.7144
$ws.Range("I30").Value = 4341
$ws.Range("J30").Value = 9999.5
$ws.Range("K30").Value = 4341
$ws.Range("L30").Value = 9999.5
$ws.Range("M30").Value = -4233
$ws.Range("N30").Value = -10215.5
# Row 68
$ws.Range("H68").Value = 6409.7617
$ws.Range("I68").Value = 4933.8335
$ws.Range("J68").Value = 8377.666999999999
$ws.Range("K68").Value = 4933.8335
$ws.Range("L68").Value = 8377.666999999999
$ws.Range("M68").Value = -4184.8335
$ws.Range("N68").Value = -9875.666999999999
# Row 71
$ws.Range("H71").Value = 6409.7617
$ws.Range("I71").Value = 4933.8335
$ws.Range("J71").Value = 8377.666999999999
$ws.Range("K71").Value = 24669.1675
$ws.Range("L71").Value = 41888.335
$ws.Range("M71").Value = -20925.1675
$ws.Range("N71").Value = -49376.335
# Row 82
$ws.Range("H82").Value = 1368.6923
$ws.Range("I82").Value = 835.875
$ws.Range("J82").Value = 2221.2
$ws.Range("K82").Value = 835.875
$ws.Range("L82").Value = 2221.2
$ws.Range("M82").Value = -474.875
$ws.Range("N82").Value = -2943.2
# Row 85
$ws.Range("H85").Value = 1368.6923
$ws.Range("I85").Value = 835.875
$ws.Range("J85").Value = 2221.2
$ws.Range("K85").Value = 835.875
$ws.Range("L85").Value = 2221.2
$ws.Range("M85").Value = 412.125
$ws.Range("N85").Value = -4717.2
# Row 93
$ws.Range("H93").Value = 26007200
$ws.Range("I93").Value = 6715.143
$ws.Range("K93").Value = 6715.143
$ws.Range("M93").Value = -5467.143
# Row 136
$ws.Range("H136").Value = 5788.8
$ws.Range("I136").Value = 3269
$ws.Range("K136").Value = 9807
$ws.Range("M136").Value = -7257

$ws = $wb.Worksheets.Item("WVR")
# Row 49
$ws.Range("H49").Value = 27746.084
$ws.Range("I49").Value = 24998
$ws.Range("J49").Value = 27995.908
$ws.Range("K49").Value = 24998
$ws.Range("L49").Value = 27995.908
$ws.Range("M49").Value = -24768
$ws.Range("N49").Value = -28455.908
# Row 113
$ws.Range("H113").Value = 756.8461
$ws.Range("I113").Value = 744.6
$ws.Range("K113").Value = 2233.8
$ws.Range("M113").Value = -63.80000000000018
# Row 122
$ws.Range("H122").Value = 4218.8125
$ws.Range("I122").Value = 2708.4167
$ws.Range("J122").Value = 8750
$ws.Range("K122").Value = 8125.250100000001
$ws.Range("L122").Value = 26250
$ws.Range("M122").Value = -5675.250100000001
$ws.Range("N122").Value = -31150
# Row 126
$ws.Range("H126").Value = 1441.7778
$ws.Range("J126").Value = 1765.3334
$ws.Range("L126").Value = 5296.0002
$ws.Range("N126").Value = -10236.0002

Write-Host "Applied price updates across all sheets."
